$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Total" row marks: Right total and Max (Correct/Total) text
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 80
$ws.Range("E12").Value = "80/140"
